# Swap the "Valor Mora" values between the 2103 period row (16) and the
# 2003 period row (28) on Hoja1, as part of updating the EC database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

$valF16 = $ws.Range("F16").Value2
$valF28 = $ws.Range("F28").Value2

$ws.Range("F16").Value2 = $valF28
$ws.Range("F28").Value2 = $valF16
